# Potturi_LabExam03Grading - fill in "Points for grading" (column E) values
# for the "Customer Class" and "Product Class" rubric sections, matching the
# "Total Points" (column D) already entered, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Customer Class section (rows 3-6)
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 2

# Product Class section (rows 10-14)
$ws.Range("E10").Value = 2
$ws.Range("E11").Value = 2
$ws.Range("E12").Value = 2
$ws.Range("E13").Value = 2
$ws.Range("E14").Value = 2

# Move the cursor / selection to E15 and scroll back to the top of the sheet
$ws.Activate()
$ws.Range("E15").Select()
$excel.ActiveWindow.ScrollRow = 1
